# Update "dados/usuarios.xlsx" sheet:
#  - C2: text "123" -> numeric 123
#  - New row 3: A3=2, B3="Rhubi", C3="123" (text), D3="klsp.201409@gmail.com"
#  - dimension grows from A1:D2 to A1:D3 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: was stored as text "123"; now becomes a real numeric value.
$ws.Range("C2").Value = 123

# New row 3.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Rhubi"

# C3 must stay a text value "123" (not become numeric like C2).
# Force text formatting so Excel doesn't auto-convert the digit string to a
# number, then restore the default "Normal" style so no extra formatting is
# left behind on the cell.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "123"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "klsp.201409@gmail.com"
